# Remove the trailing "blank / page-break / copyright notice" block that
# used to sit right after the last bibliography entry.
#
# Paragraph sequence (before):
#   ... "Prentice Hall, 2009. ... Brasiliense, 2001."   <- keep
#   (empty paragraph)                                    <- remove
#   (empty paragraph, pageBreakBefore)                   <- remove
#   "© 2020 . Contact: luizeleno@usp.br. ..."             <- remove
#   (empty paragraph)                                    <- keep
#   (empty paragraph, pageBreakBefore)                   <- keep
#
$d = $word.ActiveDocument

$bibliText = "Prentice Hall, 2009. TOLEDO, F. O que são Recursos Humanos- Primeiros Passos. São Paulo: Brasiliense, 2001."
$copyText  = "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

# Locate the two anchor paragraphs by their exact text so the edit is robust
# to the surrounding paragraph count/indices.
$n = $d.Paragraphs.Count
$bibliIdx = -1
$copyIdx = -1
for ($i = 1; $i -le $n; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $bibliText) { $bibliIdx = $i }
    if ($t -eq $copyText)  { $copyIdx = $i }
}

if ($bibliIdx -ne -1 -and $copyIdx -ne -1) {
    # Delete everything from right after the bibliography paragraph's own
    # paragraph mark through to (and including) the copyright paragraph's
    # paragraph mark -- this removes the two blank/page-break paragraphs in
    # between as well as the copyright paragraph itself, while leaving the
    # bibliography paragraph and everything after the copyright block intact.
    $startPos = $d.Paragraphs.Item($bibliIdx).Range.End
    $endPos = $d.Paragraphs.Item($copyIdx).Range.End
    $range = $d.Range($startPos, $endPos)
    $range.Delete()
}
